$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as scraped on 2023-08-10
$ws.Range("D2").Value = "29.423.82"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.848.18"
$ws.Range("E3").Value = "  -0.04%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.06%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.79"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.94%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6318"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -3.27%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -0.05%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07600"
$ws.Range("D8").Style = $origStyle
$ws.Range("E9").Value = "  -0.26%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.50"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "2.391.24"
$ws.Range("E11").Value = "  +28.73%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07724"
$ws.Range("D12").Style = $origStyle
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.989"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -0.58%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6865"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "2.480.88"
$ws.Range("E15").Value = "  +15.72%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.99"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -0.78%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009903"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +4.38%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.167"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "29.438.28"
$ws.Range("E19").Value = "  -0.25%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.82"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -2.40%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.52"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -0.61%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -0.09%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.616"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -1.20%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -0.01%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.39"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  -1.92%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.469"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  -0.58%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.472"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.75%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05819"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -4.08%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.256"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.96%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.128"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -0.08%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.026"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -1.14%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.867"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.71%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.161"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -1.90%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7179"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -0.92%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.590"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "1.247.21"
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("B39").Value = "RocketPoolETH"
$ws.Range("C39").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D39").Value = "2.428.65"
$ws.Range("E39").Value = "  +18.60%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.795"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01804"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +1.17%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9070"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +0.19%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.111"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -2.20%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.05%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.29"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +1.19%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.47"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -0.35%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.320"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000121"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.195"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4017"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.699"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +2.38%  "
